$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (preserves the cell's original
# number format / style instead of letting Excel auto-convert a
# percent-looking string into a numeric percentage).
function Set-TextValue($ws, $addr, $text) {
    $buf = $ws.Range("Z100")
    $ws.Range($addr).Copy()
    $buf.PasteSpecial(-4122)
    $ws.Range($addr).Value = "'" + $text
    $buf.Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $buf.Clear()
}

# --- Row 2: reorder "Recorded By" list ---
$ws.Range("G2").Value = "gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, System"

# --- Row 3 ---
$ws.Range("G3").Value = "majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, System"

# --- Row 4 ---
$ws.Range("G4").Value = "majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"

# --- Row 5 ---
$ws.Range("G5").Value = "asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"

# --- Row 6 ---
$ws.Range("G6").Value = "majorelle.magdy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm"
$ws.Range("L6").Value = 25

# --- Row 7 ---
$ws.Range("G7").Value = "AbeerRagheb@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg"

# --- Row 8 ---
$ws.Range("L8").Value = 2

# --- Row 9 / 10: coverage + avg attendance percentages (kept as text) ---
Set-TextValue $ws "L9" "86.2%"
Set-TextValue $ws "L10" "26.7%"

# --- Row 12 ---
$ws.Range("G12").Value = "dina.adel@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg"

# --- Row 15 (PARASITOLOGY summary row, mirrors rows 6-10 stats) ---
$ws.Range("G15").Value = "mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
$ws.Range("O15").Value = 25
$ws.Range("Q15").Value = 2
Set-TextValue $ws "R15" "86.2%"
Set-TextValue $ws "S15" "26.7%"

# --- Row 20: session got recorded (style Pending->Recorded, fill + values) ---
$ws.Range("A6:I6").Copy()
$ws.Range("A20:I20").PasteSpecial(-4122)
$ws.Range("G20").Value = "mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
$ws.Range("H20").Value = "25/251"
$ws.Range("I20").Value = "Recorded"

# --- Row 25 ---
$ws.Range("G25").Value = "Noran.Mahmoud@med.asu.edu.eg, menna-allah.gamil@med.asu.edu.eg"

# --- Row 27 ---
$ws.Range("G27").Value = "hana.amr@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"

# --- Row 30 ---
$ws.Range("G30").Value = "shorokmohamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"
